$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 odds updates (FlashScore 2024-10-17 weekly games file)
$updates = @{
    "G7"  = 2.95
    "I7"  = 2.37
    "J7"  = 3.5
    "L7"  = 2.85
    "M7"  = 1.07
    "N7"  = 7.85
    "P7"  = 2.72
    "U7"  = 1.75
    "V7"  = 1.85
    "W7"  = 8
    "Y7"  = 10.75
    "Z7"  = 37
    "AA7" = 28
    "AB7" = 37
    "AC7" = 8.25
    "AH7" = 7.5
    "AJ7" = 9.25
    "AK7" = 25
    "AL7" = 20
    "AM7" = 30
    "AN7" = 4.85
    "AO7" = 16
    "AP7" = 22
    "AQ7" = 75
    "AR7" = 110
    "AT7" = 2.57
    "AW7" = 4.25
    "AX7" = 11.75
    "AY7" = 18
    "AZ7" = 45
    "BA7" = 70
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
